$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.413.24'
$ws.Range("E2").Value = '  +3.72%  '

$ws.Range("D3").Value = '2.258.25'
$ws.Range("E3").Value = '  +1.51%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.73%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.517'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.479'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '31.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +14.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0804'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.108'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.40%  '

$ws.Range("D15").Value = '2.603.13'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.39'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.80%  '

$ws.Range("D17").Value = '2.257.12'
$ws.Range("E17").Value = '  +0.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.743'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.31%  '

$ws.Range("D19").Value = '40.330.66'
$ws.Range("E19").Value = '  +3.66%  '

$ws.Range("D20").Value = '0.0₃0897'
$ws.Range("E20").Value = '  +4.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.99%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +10.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.49'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.97%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.98'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.95'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0724'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.40'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +16.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.103'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.112'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.38%  '

$ws.Range("D43").Value = '2.020.72'
$ws.Range("E43").Value = '  +6.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0274'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +12.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.09%  '

$ws.Range("D49").Value = '2.471.82'
$ws.Range("E49").Value = '  +1.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +16.34%  '
